# Day 7 (Tuesday of week 2) content fill-in.
#
# Slides 31-36 (p:sldId 271-276, layout "1_Title Slide") were blank
# placeholder/"template" day-title slides (empty title, author
# "Peter Mackenzie-Helnwein" / "University of Washington" left over from
# copying an earlier day's slide). This pass gives them their real
# Day 7 titles and switches the presenter block over to Frank McKenna /
# UC Berkeley, matching the other Frank McKenna sessions already in the
# deck.

$p = $ppt.ActivePresentation

function Set-DayTitleSlide($SlideIndex, $TitleRuns, $Author, $Org) {
    $slide = $p.Slides.Item($SlideIndex)

    # Shape 1: Subtitle 6 (the big slide title placeholder)
    $titleShape = $slide.Shapes.Item(1)
    $titleTr = $titleShape.TextFrame.TextRange
    if ($TitleRuns.Length -gt 0) {
        $full = [string]::Join("", $TitleRuns)
        $titleTr.Text = $full
        if ($TitleRuns.Length -gt 1) {
            # Split the combined text back into separate runs, mirroring
            # how PowerPoint keeps a differently-flagged word (e.g. a
            # spell-check exception like "OpenSees") in its own run.
            $pos = 1
            foreach ($chunk in $TitleRuns) {
                if ($chunk.Length -gt 0) {
                    $sub = $titleTr.Characters($pos, $chunk.Length)
                    $sub.Text = $chunk
                }
                $pos = $pos + $chunk.Length
            }
        }
    }

    # Shape 2: Text Placeholder 4 (presenter name)
    $authorShape = $slide.Shapes.Item(2)
    if ([string]::IsNullOrEmpty($Author)) {
        $authorShape.TextFrame.TextRange.Text = ""
    } else {
        $authorShape.TextFrame.TextRange.Text = $Author
    }

    # Shape 3: Text Placeholder 5 (presenter organization)
    $orgShape = $slide.Shapes.Item(3)
    $orgShape.TextFrame.TextRange.Text = $Org
}

# Slide 31 - "Abstraction in C"
Set-DayTitleSlide 31 @("Abstraction in C") "Frank McKenna" "University of California at Berkeley"

# Slide 32 - "Object Oriented Programming in C++"
Set-DayTitleSlide 32 @("Object Oriented Programming in C++") "Frank McKenna" "University of California at Berkeley"

# Slide 33 - "Software Design Example: OpenSees"
Set-DayTitleSlide 33 @("Software Design Example: ", "OpenSees") "Frank McKenna" "University of California at Berkeley"

# Slide 34 - "EMACS - Tips & Tricks"
Set-DayTitleSlide 34 @("EMACS - Tips & Tricks") "Frank McKenna" "University of California at Berkeley"

# Slide 35 - "Coding Exercise: a Vector class"
Set-DayTitleSlide 35 @("Coding Exercise: a ", "Vector class") "Frank McKenna" "University of California at Berkeley"

# Slide 36 - divider slide: no title text, author cleared, org updated.
Set-DayTitleSlide 36 @() "" "University of California at Berkeley"
